$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Travail")

# Fill in row 14 with the new "Rendu 4" entry (left side: Gabriel, right side: Stan)
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

$ws.Range("B14").Value = 45412
$ws.Range("C14").Value = "25m"
$ws.Range("D14").Value = "modifications dans le graph d'état"

$ws.Range("F14").Value = 45412
$ws.Range("G14").Value = "25m"
$ws.Range("H14").Value = "modifications dans le graph d'état"

# Update selection to H14
$ws.Range("H14").Select()
